$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("27").Copy($ws.Rows("28"))
